# Slide 7 (1-based) is missing the "Slide Number Placeholder" shape that
# every other content slide in this deck already has. Re-add it the same
# way PowerPoint does when you turn on "Slide Number" in
# Insert > Header & Footer for this slide: via HeadersFooters.SlideNumber.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# The slide's existing shape ids are 2, 4, 6 (title, footer, content).
# PowerPoint's per-slide "next shape id" allocator here hands out the
# smallest free id >= 2 the first few times a shape is created, which
# would land on 3 rather than the 5 used in the authored file. Add and
# immediately discard a throwaway shape to consume id 3, so the
# placeholder we actually want lands on id 5 (matching "Slide Number
# Placeholder 4" / id="5" in the target deck).
$dummy = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$dummy.Delete()

# Turn on the slide-number field for this slide; this creates the
# sldNum placeholder shape (type="sldNum" sz="quarter" idx="12") with
# the dynamic <a:fld type="slidenum"> field inside it.
$hf = $s.HeadersFooters
$sn = $hf.SlideNumber
$sn.Visible = $true

# Locate the placeholder just created (ppPlaceholderSlideNumber = 13)
# rather than assuming it is the last shape in the collection.
$ph = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Type -eq 14 -and $candidate.PlaceholderFormat.Type -eq 13) {
        $ph = $candidate
    }
}
if ($ph -eq $null) {
    $ph = $s.Shapes.Item($s.Shapes.Count)
}
$ph.Name = "Slide Number Placeholder 4"

# Match the explicit position/size from the authored slide
# (6553200,4705350) / (2133600,357188) EMU -> points (1 pt = 12700 EMU).
$ph.Left = 516
$ph.Top = 370.5
$ph.Width = 168
$ph.Height = 28.12503937007874
